# Update countries & provincias Spain
# Daily refresh of the "Pais" (countries) COVID dashboard sheet:
#  - A handful of countries swapped ranking position (each row's rank is
#    fixed, so when two countries cross in ranking the country name shown
#    on a row changes while its numeric columns follow the country).
#  - Updated totals (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes) for the affected rows.
#  - Updated "last refreshed" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Country, $B, $C, $D, $E, $G, $H) {
    if ($Country -ne $null) { $ws.Cells.Item($Row, 1).Value = $Country }
    if ($B -ne $null) { $ws.Cells.Item($Row, 2).Value = $B }
    if ($C -ne $null) { $ws.Cells.Item($Row, 3).Value = $C }
    if ($D -ne $null) { $ws.Cells.Item($Row, 4).Value = $D }
    if ($E -ne $null) { $ws.Cells.Item($Row, 5).Value = $E }
    if ($G -ne $null) { $ws.Cells.Item($Row, 7).Value = $G }
    if ($H -ne $null) { $ws.Cells.Item($Row, 8).Value = $H }
}

# --- Rows whose underlying country changed rank (name swap) plus new numbers ---

# Etiopia <-> Costa Rica
Set-Row 54 "Etiopia"    77860 872   32325  44321        6   1214
Set-Row 55 "Costa Rica" 77829 $null 42621  34278        $null 930

# Camerun / Costa de Marfil / Tunez three-way rotation
Set-Row 84 "Tunez"           20944 1223  5032   15636   5     276
Set-Row 85 "Camerun"         20924 $null 19764  740     $null 420
Set-Row 86 "Costa de Marfil" 19849 56    19421  308     $null 120

# Guinea <-> Consejo Danes para los Refugiados
Set-Row 102 "Consejo Danes para los Refugiados" 10752 23    10212 266 2     274
Set-Row 103 "Guinea"                            10735 $null 10066 603 $null 66

# Siria <-> Bahamas
Set-Row 134 "Bahamas" 4332 112   2375 1861 $null 96
Set-Row 135 "Siria"   4329 40    1143 2982 1     204

# Nueva Caledonia <-> Santa Lucia (numbers unchanged, only rank/name swap)
Set-Row 207 "Santa Lucia"    $null $null $null $null $null $null
Set-Row 208 "Nueva Caledonia" $null $null $null $null $null $null

# --- Rows with only numeric updates (no country-name / rank change) ---

Set-Row 4   $null 7588513 39190 4807847 2566619 525   214047
Set-Row 5   $null 6547413 75479 5506732 938869  937   101812
Set-Row 13  $null 679716  1883  612763  50015   29    16938
Set-Row 14  $null $null   $null $null   476649  49    32198
Set-Row 25  $null 299787  1424  $null   30690   $null $null
Set-Row 27  $null 264443  5523  191251  71510   49    1682
Set-Row 119 $null $null   $null 4539    1065    $null $null
Set-Row 131 $null 4852    5     3211    1612    $null $null
Set-Row 140 $null 3590    5     2226    1249    $null $null
Set-Row 161 $null 1840    22    1375    417     $null $null
Set-Row 166 $null 1214    3     1075    54      $null $null
Set-Row 185 $null 342     1     $null   3       $null $null

# --- Last-updated timestamp (shared-string cell A1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Octubre de 2020 a las 22:23"
